$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.099562666666667
$ws.Range("H2").Value = 9.298688
$ws.Range("I2").Value = 0.2686390288432488
$ws.Range("J2").Value = 0.2686390288432488
$ws.Range("M2").Value = 8.813278666666667
$ws.Range("N2").Value = 26.439836
$ws.Range("O2").Value = 0.3770976991891536
$ws.Range("P2").Value = 0.3770976991891536
$ws.Range("Q2").Value = 27.31730952612978
$ws.Range("R2").Value = 245.855785735168
$ws.Range("S2").Value = 0.1013031596891978
$ws.Range("T2").Value = 0.1013031596891978
$ws.Range("G3").Value = 3.099562666666667
$ws.Range("H3").Value = 9.298688
$ws.Range("I3").Value = 0.2686390288432488
$ws.Range("J3").Value = 0.2686390288432488
$ws.Range("O3").Value = 0.5522024902836482
$ws.Range("P3").Value = 0.5522024902836482
$ws.Range("Q3").Value = 40.00206413513956
$ws.Range("R3").Value = 360.018577216256
$ws.Range("S3").Value = 0.1483431407146228
$ws.Range("T3").Value = 0.1483431407146228
$ws.Range("G4").Value = 3.099562666666667
$ws.Range("H4").Value = 9.298688
$ws.Range("I4").Value = 0.2686390288432488
$ws.Range("J4").Value = 0.2686390288432488
$ws.Range("M4").Value = 1.649921333333333
$ws.Range("N4").Value = 4.949764
$ws.Range("O4").Value = 0.07059592260441032
$ws.Range("P4").Value = 0.07059592260441033
$ws.Range("Q4").Value = 5.11403456773689
$ws.Range("R4").Value = 46.026311109632
$ws.Range("S4").Value = 0.01896482008874194
$ws.Range("T4").Value = 0.01896482008874194
$ws.Range("G5").Value = 3.099562666666667
$ws.Range("H5").Value = 9.298688
$ws.Range("I5").Value = 0.2686390288432488
$ws.Range("J5").Value = 0.2686390288432488
$ws.Range("M5").Value = 0.002428
$ws.Range("N5").Value = 0.007284
$ws.Range("O5").Value = 0.0001038879227879399
$ws.Range("P5").Value = 0.0001038879227879399
$ws.Range("Q5").Value = 0.007525738154666667
$ws.Range("R5").Value = 0.067731643392
$ws.Range("S5").Value = 0.0000279083506862946
$ws.Range("T5").Value = 0.0000279083506862946
$ws.Range("G6").Value = 6.189892666666666
$ws.Range("I6").Value = 0.5364778626674904
$ws.Range("J6").Value = 0.5364778626674905
$ws.Range("M6").Value = 8.813278666666667
$ws.Range("N6").Value = 26.439836
$ws.Range("O6").Value = 0.3770976991891536
$ws.Range("P6").Value = 0.3770976991891536
$ws.Range("Q6").Value = 54.55324898808978
$ws.Range("R6").Value = 490.979240892808
$ws.Range("S6").Value = 0.2023045676778254
$ws.Range("T6").Value = 0.2023045676778254
$ws.Range("G7").Value = 6.189892666666666
$ws.Range("I7").Value = 0.5364778626674904
$ws.Range("J7").Value = 0.5364778626674905
$ws.Range("O7").Value = 0.5522024902836482
$ws.Range("P7").Value = 0.5522024902836482
$ws.Range("Q7").Value = 79.88497413020956
$ws.Range("R7").Value = 718.964767171886
$ws.Range("S7").Value = 0.2962444117470373
$ws.Range("T7").Value = 0.2962444117470373
$ws.Range("G8").Value = 6.189892666666666
$ws.Range("I8").Value = 0.5364778626674904
$ws.Range("J8").Value = 0.5364778626674905
$ws.Range("M8").Value = 1.649921333333333
$ws.Range("N8").Value = 4.949764
$ws.Range("O8").Value = 0.07059592260441032
$ws.Range("P8").Value = 0.07059592260441033
$ws.Range("Q8").Value = 10.21283596177689
$ws.Range("R8").Value = 91.915523655992
$ws.Range("S8").Value = 0.03787314967185362
$ws.Range("T8").Value = 0.03787314967185364
$ws.Range("G9").Value = 6.189892666666666
$ws.Range("I9").Value = 0.5364778626674904
$ws.Range("J9").Value = 0.5364778626674905
$ws.Range("M9").Value = 0.002428
$ws.Range("N9").Value = 0.007284
$ws.Range("O9").Value = 0.0001038879227879399
$ws.Range("P9").Value = 0.0001038879227879399
$ws.Range("Q9").Value = 0.01502905939466667
$ws.Range("R9").Value = 0.135261534552
$ws.Range("S9").Value = 0.0000557335707742393
$ws.Range("T9").Value = 0.00005573357077423931
$ws.Range("G10").Value = 1.888584
$ws.Range("H10").Value = 5.665752
$ws.Range("I10").Value = 0.1636835341659699
$ws.Range("J10").Value = 0.1636835341659699
$ws.Range("M10").Value = 8.813278666666667
$ws.Range("N10").Value = 26.439836
$ws.Range("O10").Value = 0.3770976991891536
$ws.Range("P10").Value = 0.3770976991891536
$ws.Range("Q10").Value = 16.644617077408
$ws.Range("R10").Value = 149.801553696672
$ws.Range("S10").Value = 0.06172468412913647
$ws.Range("T10").Value = 0.06172468412913648
$ws.Range("G11").Value = 1.888584
$ws.Range("H11").Value = 5.665752
$ws.Range("I11").Value = 0.1636835341659699
$ws.Range("J11").Value = 0.1636835341659699
$ws.Range("O11").Value = 0.5522024902836482
$ws.Range("P11").Value = 0.5522024902836482
$ws.Range("Q11").Value = 24.373521821336
$ws.Range("R11").Value = 219.361696392024
$ws.Range("S11").Value = 0.09038645518487719
$ws.Range("T11").Value = 0.0903864551848772
$ws.Range("G12").Value = 1.888584
$ws.Range("H12").Value = 5.665752
$ws.Range("I12").Value = 0.1636835341659699
$ws.Range("J12").Value = 0.1636835341659699
$ws.Range("M12").Value = 1.649921333333333
$ws.Range("N12").Value = 4.949764
$ws.Range("O12").Value = 0.07059592260441032
$ws.Range("P12").Value = 0.07059592260441033
$ws.Range("Q12").Value = 3.116015031392
$ws.Range("R12").Value = 28.044135282528
$ws.Range("S12").Value = 0.01155539010959716
$ws.Range("T12").Value = 0.01155539010959717
$ws.Range("G13").Value = 1.888584
$ws.Range("H13").Value = 5.665752
$ws.Range("I13").Value = 0.1636835341659699
$ws.Range("J13").Value = 0.1636835341659699
$ws.Range("M13").Value = 0.002428
$ws.Range("N13").Value = 0.007284
$ws.Range("O13").Value = 0.0001038879227879399
$ws.Range("P13").Value = 0.0001038879227879399
$ws.Range("Q13").Value = 0.004585481952
$ws.Range("R13").Value = 0.041269337568
$ws.Range("S13").Value = 0.00001700474235909141
$ws.Range("T13").Value = 0.00001700474235909141
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.3599813333333333
$ws.Range("H14").Value = 1.079944
$ws.Range("I14").Value = 0.03119957432329092
$ws.Range("J14").Value = 0.03119957432329093
$ws.Range("M14").Value = 8.813278666666667
$ws.Range("N14").Value = 26.439836
$ws.Range("O14").Value = 0.3770976991891536
$ws.Range("P14").Value = 0.3770976991891536
$ws.Range("Q14").Value = 3.172615805464889
$ws.Range("R14").Value = 28.553542249184
$ws.Range("S14").Value = 0.011765287692994
$ws.Range("T14").Value = 0.011765287692994
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.3599813333333333
$ws.Range("H15").Value = 1.079944
$ws.Range("I15").Value = 0.03119957432329092
$ws.Range("J15").Value = 0.03119957432329093
$ws.Range("O15").Value = 0.5522024902836482
$ws.Range("P15").Value = 0.5522024902836482
$ws.Range("Q15").Value = 4.645815533369777
$ws.Range("R15").Value = 41.812339800328
$ws.Range("S15").Value = 0.01722848263711102
$ws.Range("T15").Value = 0.01722848263711102
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.3599813333333333
$ws.Range("H16").Value = 1.079944
$ws.Range("I16").Value = 0.03119957432329092
$ws.Range("J16").Value = 0.03119957432329093
$ws.Range("M16").Value = 1.649921333333333
$ws.Range("N16").Value = 4.949764
$ws.Range("O16").Value = 0.07059592260441032
$ws.Range("P16").Value = 0.07059592260441033
$ws.Range("Q16").Value = 0.5939408814684445
$ws.Range("R16").Value = 5.345467933216
$ws.Range("S16").Value = 0.002202562734217593
$ws.Range("T16").Value = 0.002202562734217594
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.3599813333333333
$ws.Range("H17").Value = 1.079944
$ws.Range("I17").Value = 0.03119957432329092
$ws.Range("J17").Value = 0.03119957432329093
$ws.Range("M17").Value = 0.002428
$ws.Range("N17").Value = 0.007284
$ws.Range("O17").Value = 0.0001038879227879399
$ws.Range("P17").Value = 0.0001038879227879399
$ws.Range("Q17").Value = 0.0008740346773333333
$ws.Range("R17").Value = 0.007866312096
$ws.Range("S17").Value = 0.0001033363845498802
